$wb = $excel.ActiveWorkbook
$wsDb = $wb.Worksheets.Item("DataBase")
$wsSp = $wb.Worksheets.Item("Specialist")

# --- DataBase sheet: extend "Data Table:  Examination type" ---
# Add the 7 new data rows (97-103) mirroring the existing header/row96 pattern.
$wsDb.Range("A97").Value = "physiotherapy"
$wsDb.Range("B99").Value = "Otorhinolaryngology"
$wsDb.Range("A99").Value = "sinusitis"
$wsDb.Range("A100").Value = "Hearing test"
$wsDb.Range("A101").Value = "orthodontics"
$wsDb.Range("A102").Value = "psoriasis"
$wsDb.Range("A103").Value = "Tova"
$wsDb.Range("A98").Value = "Papilloma"

$wsDb.Range("B97").Value = "Othopedy"
$wsDb.Range("B98").Value = "Gynecology"
$wsDb.Range("B100").Value = "Audiology"
$wsDb.Range("B101").Value = "Dentist"
$wsDb.Range("B102").Value = "Dermatology"
$wsDb.Range("B103").Value = "Neurology"

$wsDb.Activate()
$wsDb.Range("C100").Select()
$excel.ActiveWindow.Zoom = 189

$wsSp.Activate()
$wsSp.Range("D25").Select()

$wsDb.Activate()
